# Changes to read the password from excelsheet
# - Populate the SignIn sheet's row 4 with the login values used by the
#   test (username / password / confirmation message) and make SignIn the
#   active (selected) sheet/tab, matching the new workbook state.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("SignIn")

# Fill in the username / password / login message that the automated
# tests now read from the sheet instead of having blank cells.
$ws.Range("A4").Value = "maya"
$ws.Range("A4").Font.ThemeColor = 1

$ws.Range("B4").Value = "babuji_123"
$ws.Range("B4").Font.ThemeColor = 1

$ws.Range("C4").Value = "You are looged in"
$ws.Range("C4").Font.ThemeColor = 1

# Make the SignIn sheet the active tab/selection (previously Register was
# the active tab) and select C4, which is where the cursor ends up after
# entering the new data.
$ws.Activate()
$ws.Range("C4").Select()
